# Adapt column header formatting to respective input-file names:
#   "<field>_old" -> "<field>_FV2304"
#   "<field>_new" -> "<field>_FV2310"
# Then wrap the used range in a named table and freeze the header row,
# matching the regenerated AHB-diff xlsx export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-suffix the two header blocks (A1:J1 "_old" -> "_FV2304",
#        L1:U1 "_new" -> "_FV2310"; K1 "diff" is untouched) -------------
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2304Headers[$i]
}
for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2310Headers[$i]
}

# --- 2) Turn the used range into "Table1" (adds xl/tables/table1.xml,
#        the sheet's tableParts entry and its part relationship) -------
$usedRange = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $usedRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3) Freeze the header row (pane ySplit=1, state=frozen) -----------
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
